$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ G=6.240107999999999;  H=18.720324;  I=0.01732230523539376; J=0.01732230523539376;
             M=11.146846; N=33.440538; O=0.2594806085672136; P=0.2594806085672136;
             Q=69.55752289936801; R=626.017706094312; S=0.004494802304267002; T=0.004494802304267003 }
    3  = @{ G=6.240107999999999;  H=18.720324;  I=0.01732230523539376; J=0.01732230523539376;
             O=0.6444737471070977; P=0.6444737471070977;
             Q=172.760491313676; R=1554.844421823084; S=0.01116377096358711; T=0.01116377096358711 }
    4  = @{ G=6.240107999999999;  H=18.720324;  I=0.01732230523539376; J=0.01732230523539376;
             O=0.09604564432568881; P=0.09604564432568881;
             Q=25.746421443428; R=231.717792990852; S=0.001663731967539646; T=0.001663731967539647 }
    5  = @{ I=0.9592798330716089; J=0.9592798330716091;
             M=11.146846; N=33.440538; O=0.2594806085672136; P=0.2594806085672136;
             Q=3851.977438859834; R=34667.79694973851; S=0.2489145148716761; T=0.2489145148716762 }
    6  = @{ I=0.9592798330716089; J=0.9592798330716091;
             O=0.6444737471070977; P=0.6444737471070977;
             S=0.6182306685439309; T=0.6182306685439312 }
    7  = @{ I=0.9592798330716089; J=0.9592798330716091;
             O=0.09604564432568881; P=0.09604564432568881;
             S=0.09213464965600188; T=0.0921346496560019 }
    8  = @{ G=8.428738666666666; I=0.02339786169299727; J=0.02339786169299728;
             M=11.146846; N=33.440538; O=0.2594806085672136; P=0.2594806085672136;
             Q=93.95385189157868; R=845.5846670242081; S=0.006071291391270425; T=0.006071291391270427 }
    9  = @{ G=8.428738666666666; I=0.02339786169299727; J=0.02339786169299728;
             O=0.6444737471070977; P=0.6444737471070977;
             S=0.01507930759957957; T=0.01507930759957958 }
    10 = @{ G=8.428738666666666; I=0.02339786169299727; J=0.02339786169299728;
             O=0.09604564432568881; P=0.09604564432568881;
             Q=34.77661892206311; R=312.989570298568; S=0.002247262702147275; T=0.002247262702147276 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $values[$row][$col]
    }
}
